$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.893.45"
$ws.Range("E2").Value = "  -0.10%  "
$ws.Range("D3").Value = "2.363.75"
$ws.Range("E3").Value = "  +0.50%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'0.673"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.73%  "
$ws.Range("D6").Value = "'240.42"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.71%  "
$ws.Range("D7").Value = "'74.25"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.69%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'0.616"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.40%  "
$ws.Range("D10").Value = "'0.103"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.56%  "
$ws.Range("D11").Value = "'60.78"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.12%  "
$ws.Range("D12").Value = "'37.81"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +16.92%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'7.32"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.38%  "
$ws.Range("B14").Value = "TRON"
$ws.Range("C14").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D14").Value = "'0.109"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.83%  "
$ws.Range("D15").Value = "'16.38"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.15%  "
$ws.Range("D16").Value = "'0.922"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.88%  "
$ws.Range("D17").Value = "2.356.40"
$ws.Range("E17").Value = "  +0.26%  "
$ws.Range("D18").Value = "43.864.25"
$ws.Range("E18").Value = "  +0.03%  "
$ws.Range("D19").Value = "'0.0000104"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.30%  "
$ws.Range("D20").Value = "'78.05"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.70%  "
$ws.Range("D21").Value = "'6.58"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.97%  "
$ws.Range("D22").Value = "'254.13"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.41%  "
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("E24").Value = "  +3.35%  "
$ws.Range("E25").Value = "  -2.44%  "
$ws.Range("E26").Value = "  +0.91%  "
$ws.Range("D27").Value = "'10.54"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.76%  "
$ws.Range("E28").Value = "  +0.51%  "
$ws.Range("D29").Value = "'22.38"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.03%  "
$ws.Range("D30").Value = "'175.81"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.23%  "
$ws.Range("E31").Value = "  +0.81%  "
$ws.Range("E32").Value = "  -1.18%  "
$ws.Range("D33").Value = "'0.0754"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.12%  "
$ws.Range("D34").Value = "'5.10"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.98%  "
$ws.Range("D35").Value = "'5.41"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.59%  "
$ws.Range("D36").Value = "'3.82"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.89%  "
$ws.Range("D37").Value = "'6.63"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.43%  "
$ws.Range("D38").Value = "'2.42"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.03%  "
$ws.Range("E39").Value = "  +0.69%  "
$ws.Range("D40").Value = "'5.44"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +15.08%  "
$ws.Range("D41").Value = "'20.41"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.95%  "
$ws.Range("D42").Value = "'65.10"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +12.46%  "
$ws.Range("D43").Value = "'0.205"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.19%  "
$ws.Range("B44").Value = "Cronos"
$ws.Range("C44").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D44").Value = "'0.107"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.56%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "'9.09"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.84%  "
$ws.Range("D46").Value = "'2.56"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.00%  "
$ws.Range("E47").Value = "  +0.14%  "
$ws.Range("E48").Value = "  +0.21%  "
$ws.Range("E49").Value = "  -0.94%  "
$ws.Range("D50").Value = "'98.51"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.32%  "
$ws.Range("D51").Value = "'4.43"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +16.78%  "
